# Apply "diagram labeling thru 2016" edit:
# 1. Update H9, H17, H21 subcategory labels
# 2. Delete column I (is_viewed) entirely, shrinking used range to A1:H22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update subcategory text for a few rows
$ws.Range("H9").Value = "line graph(s)"
$ws.Range("H17").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H21").Value = "data collection, data analysis, data gathering diagram"

# Delete the entire "is_viewed" column (I), shifting nothing left since it's the last column
$ws.Range("I1:I22").Delete()
